$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text: update only the changed run, leave the rest of the
# rich-text string (and its per-run formatting) untouched.
# ------------------------------------------------------------------
$a8 = $ws.Range("A8")
# "Volume 30   Number  19" -> "...  20" (volume number is the last run)
$a8.Characters(21, 2).Text = "20"

$c9 = $ws.Range("C9")
# "Report Covering the Week  5/8/2023  Through  5/14/2023" -> next week's dates
$c9.Characters(27, 8).Text = "5/15/2023"
$t9 = $c9.Text
$idx9 = $t9.IndexOf("5/14/2023")
$c9.Characters($idx9 + 1, 9).Text = "5/21/2023"

# ------------------------------------------------------------------
# Cells that flip between a numeric figure and the shared "0" /
# "***.*" placeholder text. Copy a same-styled donor cell so the
# destination picks up the right style/number-format; when the new
# value differs from the donor's, stamp the real number in afterwards.
# ------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("D15"))

$ws.Range("E14").Copy($ws.Range("E15"))

$ws.Range("C16").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1

$ws.Range("H16").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 0

$ws.Range("D14").Copy($ws.Range("C20"))

$ws.Range("F14").Copy($ws.Range("C22"))

$ws.Range("I14").Copy($ws.Range("C23"))

$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1

$ws.Range("C14").Copy($ws.Range("D27"))

$ws.Range("E14").Copy($ws.Range("E27"))

$ws.Range("C14").Copy($ws.Range("F28"))

$ws.Range("C14").Copy($ws.Range("F29"))

$ws.Range("C14").Copy($ws.Range("F30"))

# ------------------------------------------------------------------
# Plain numeric refreshes (new crime-count / percent-change figures)
# ------------------------------------------------------------------
$ws.Range("M15").Value = -14.285714285714
$ws.Range("N15").Value = -40
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 29
$ws.Range("K16").Value = 62.068965517241
$ws.Range("L16").Value = 42.424242424242
$ws.Range("M16").Value = -7.843137254901
$ws.Range("N16").Value = -78.341013824884
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 46.153846153846
$ws.Range("I17").Value = 83
$ws.Range("J17").Value = 88
$ws.Range("K17").Value = -5.681818181818
$ws.Range("L17").Value = -29.661016949152
$ws.Range("M17").Value = 29.6875
$ws.Range("N17").Value = -47.798742138364
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 28
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = 7.692307692307
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = -34.883720930232
$ws.Range("N18").Value = -89.855072463768
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 100
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -29.411764705882
$ws.Range("I19").Value = 65
$ws.Range("J19").Value = 78
$ws.Range("K19").Value = -16.666666666666
$ws.Range("L19").Value = 20.37037037037
$ws.Range("M19").Value = 75.675675675675
$ws.Range("N19").Value = -41.441441441441
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("L20").Value = 19.047619047619
$ws.Range("M20").Value = -19.354838709677
$ws.Range("N20").Value = -88.151658767772
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 48
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = 9.090909090909
$ws.Range("I21").Value = 254
$ws.Range("J21").Value = 246
$ws.Range("K21").Value = 3.252032520325
$ws.Range("L21").Value = -1.550387596899
$ws.Range("M21").Value = 7.17299578059
$ws.Range("N21").Value = -74.343434343434
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -20
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = 26.923076923076
$ws.Range("L23").Value = -15.384615384615
$ws.Range("M23").Value = 94.117647058823
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 46.153846153846
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 3.508771929824
$ws.Range("I24").Value = 219
$ws.Range("J24").Value = 206
$ws.Range("K24").Value = 6.31067961165
$ws.Range("L24").Value = 5.288461538461
$ws.Range("M24").Value = 64.661654135338
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 83.333333333333
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 146
$ws.Range("J25").Value = 140
$ws.Range("K25").Value = 4.285714285714
$ws.Range("L25").Value = -16.571428571428
$ws.Range("M25").Value = -5.194805194805
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 15
$ws.Range("L26").Value = 15.384615384615
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 15
$ws.Range("K27").Value = -6.25
$ws.Range("L27").Value = -21.052631578947
$ws.Range("L28").Value = -25
$ws.Range("L29").Value = -50
